# Applies the change described by the diff:
#  - A new parameter column "param_E_pv3_solar" is inserted logically before
#    the old column U (param_P_to_charging_station1). The two existing
#    headers shift right by one column (old U -> V, old V -> W, a brand new
#    column W), and the sheet grows from 22 to 23 columns (A:V -> A:W).
#  - The data rows (2-17) are updated with the recomputed values for the
#    (new) U, V and W columns from the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column is inserted logically before the old "U" column (which held
# param_P_to_charging_station1). Shift the two existing headers one column to
# the right (U->V, V->W) and put the brand-new header in U.
$oldU1 = $ws.Range("U1").Value2   # param_P_to_charging_station1
$oldV1 = $ws.Range("V1").Value2   # param_P_to_charging_station2

$ws.Range("W1").Value = $oldV1
$ws.Range("V1").Value = $oldU1
$ws.Range("U1").Value = "param_E_pv3_solar"     # new header

# Copy the header's formatting (bold/centered/bordered) onto the new W1 cell
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)  # xlPasteFormats

# Row -> (U, V, W) values
$rowValues = @{
    2  = @(0.12, 0, 109.395)
    3  = @(0.12, 59.89148611111113, 49.33500000000002)
    4  = @(0.12, 54.11648611111113, 0)
    5  = @(0.12, 0, 0)
    6  = @(0.12, 71.51084722222224, 0)
    7  = @(0.12, 30.25000000000002, 0)
    8  = @(0.12, 0, 0)
    9  = @(0.12, 0, 0)
    10 = @(0.12, 0, 0)
    11 = @(0.12, 0, 0)
    12 = @(0.12, 0, 0)
    13 = @(0.12, 0, 0)
    14 = @(0.12, 0, 0)
    15 = @(0.12, 0, 0)
    16 = @(0.12, 0, 0)
    17 = @(0.12, 0, 0)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    $ws.Cells.Item($r, 21).Value = $vals[0]   # column U
    $ws.Cells.Item($r, 22).Value = $vals[1]   # column V
    $ws.Cells.Item($r, 23).Value = $vals[2]   # column W
}
